$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5750
$ws.Range("J64").Value = 7000
$ws.Range("L64").Value = 7000
$ws.Range("N64").Value = -7496

$ws.Range("H67").Value = 5750
$ws.Range("J67").Value = 7000
$ws.Range("L67").Value = 7000
$ws.Range("N67").Value = -8716

$ws.Range("H86").Value = 933
$ws.Range("I86").Value = 933
$ws.Range("K86").Value = 933
$ws.Range("M86").Value = 190

$ws.Range("H89").Value = 933
$ws.Range("I89").Value = 933
$ws.Range("K89").Value = 4665
$ws.Range("M89").Value = 951

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2449.5833
$ws.Range("I32").Value = 2508.6365
$ws.Range("K32").Value = 2508.6365
$ws.Range("M32").Value = -2221.6365

$ws.Range("H61").Value = 2993.6428
$ws.Range("I61").Value = 1488.875
$ws.Range("K61").Value = 1488.875
$ws.Range("M61").Value = -1276.875

$ws.Range("H97").Value = 2119.625
$ws.Range("I97").Value = 2279.7144
$ws.Range("K97").Value = 2279.7144
$ws.Range("M97").Value = -1783.7144

$ws.Range("H136").Value = 2993.6428
$ws.Range("I136").Value = 1488.875
$ws.Range("K136").Value = 4466.625
$ws.Range("M136").Value = -1916.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1831.1428
$ws.Range("I134").Value = 1085.5294
$ws.Range("K134").Value = 3256.5882
$ws.Range("M134").Value = -721.5881999999997

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2193.4443
$ws.Range("I16").Value = 2217.75
$ws.Range("J16").Value = 1999
$ws.Range("K16").Value = 2217.75
$ws.Range("L16").Value = 1999
$ws.Range("M16").Value = -1930.75
$ws.Range("N16").Value = -2573

$ws.Range("H62").Value = 3574.5
$ws.Range("I62").Value = 3432.6667
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 3432.6667
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -2808.6667
$ws.Range("N62").Value = -5248

$ws.Range("H65").Value = 3574.5
$ws.Range("I65").Value = 3432.6667
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 17163.3335
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -14043.3335
$ws.Range("N65").Value = -26240

$ws.Range("H105").Value = 4175.3335
$ws.Range("I105").Value = 3597.6667
$ws.Range("K105").Value = 3597.6667
$ws.Range("M105").Value = -1850.6667

$ws.Range("H113").Value = 2193.4443
$ws.Range("I113").Value = 2217.75
$ws.Range("J113").Value = 1999
$ws.Range("K113").Value = 2217.75
$ws.Range("L113").Value = 1999
$ws.Range("M113").Value = -47.75
$ws.Range("N113").Value = -6339

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2630.6667
$ws.Range("I102").Value = 2310.8572
$ws.Range("K102").Value = 2310.8572
$ws.Range("M102").Value = -688.8571999999999

$ws.Range("H126").Value = 7984.25
$ws.Range("I126").Value = 7312.3335
$ws.Range("K126").Value = 21937.0005
$ws.Range("M126").Value = -19467.0005

$ws.Range("H132").Value = 3376.611
$ws.Range("I132").Value = 2913
$ws.Range("K132").Value = 8739
$ws.Range("M132").Value = -6209

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2125.5
$ws.Range("I7").Value = 2000.6666
$ws.Range("K7").Value = 2000.6666
$ws.Range("M7").Value = -1888.6666

$ws.Range("H22").Value = 1848.4166
$ws.Range("I22").Value = 1580.8334
$ws.Range("K22").Value = 1580.8334
$ws.Range("M22").Value = -1285.8334

$ws.Range("H27").Value = 1848.4166
$ws.Range("I27").Value = 1580.8334
$ws.Range("K27").Value = 1580.8334
$ws.Range("M27").Value = -1473.8334

$ws.Range("H40").Value = 5974.5
$ws.Range("I40").Value = 5974.5
$ws.Range("K40").Value = 5974.5
$ws.Range("M40").Value = -5838.5

$ws.Range("H55").Value = 1735.2
$ws.Range("I55").Value = 2157.4285
$ws.Range("J55").Value = 750
$ws.Range("K55").Value = 2157.4285
$ws.Range("L55").Value = 750
$ws.Range("M55").Value = -1984.4285
$ws.Range("N55").Value = -1096

$ws.Range("H61").Value = 2639.25
$ws.Range("I61").Value = 1499
$ws.Range("J61").Value = 3779.5
$ws.Range("K61").Value = 1499
$ws.Range("L61").Value = 3779.5
$ws.Range("M61").Value = -1297
$ws.Range("N61").Value = -4183.5

$ws.Range("H68").Value = 8110.75
$ws.Range("I68").Value = 971.5
$ws.Range("J68").Value = 15250
$ws.Range("K68").Value = 971.5
$ws.Range("L68").Value = 15250
$ws.Range("M68").Value = -222.5
$ws.Range("N68").Value = -16748

$ws.Range("H71").Value = 8110.75
$ws.Range("I71").Value = 971.5
$ws.Range("J71").Value = 15250
$ws.Range("K71").Value = 4857.5
$ws.Range("L71").Value = 76250
$ws.Range("M71").Value = -1113.5
$ws.Range("N71").Value = -83738

$ws.Range("H100").Value = 1003
$ws.Range("I100").Value = 1003
$ws.Range("K100").Value = 1003
$ws.Range("M100").Value = -462

$ws.Range("H113").Value = 2639.25
$ws.Range("I113").Value = 1499
$ws.Range("J113").Value = 3779.5
$ws.Range("K113").Value = 1499
$ws.Range("L113").Value = 3779.5
$ws.Range("M113").Value = 671
$ws.Range("N113").Value = -8119.5

$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws.Range("H126").Value = 2125.5
$ws.Range("I126").Value = 2000.6666
$ws.Range("K126").Value = 6001.9998
$ws.Range("M126").Value = -3531.9998

$ws.Range("H132").Value = 1992.8948
$ws.Range("I132").Value = 1781.8334
$ws.Range("J132").Value = 2354.7144
$ws.Range("K132").Value = 5345.5002
$ws.Range("L132").Value = 7064.1432
$ws.Range("M132").Value = -2815.5002
$ws.Range("N132").Value = -12124.1432

$ws.Range("H136").Value = 2857.6
$ws.Range("I136").Value = 2672
$ws.Range("J136").Value = 3600
$ws.Range("K136").Value = 8016
$ws.Range("L136").Value = 10800
$ws.Range("M136").Value = -5466
$ws.Range("N136").Value = -15900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 599.5714
$ws.Range("I107").Value = 532.8333
$ws.Range("K107").Value = 1598.4999
$ws.Range("M107").Value = 321.5001

$ws.Range("H113").Value = 1247
$ws.Range("I113").Value = 1331
$ws.Range("K113").Value = 3993
$ws.Range("M113").Value = -1823

$ws.Range("H136").Value = 562.2222
$ws.Range("I136").Value = 548
$ws.Range("J136").Value = 580
$ws.Range("K136").Value = 1644
$ws.Range("L136").Value = 1740
$ws.Range("M136").Value = 906
$ws.Range("N136").Value = -6840
